$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new value) pair below mirrors an updated Price/Volume(1h) figure
# from the refreshed crypto symbol list. NumberFormat is pinned to Text ("@")
# first so Excel keeps storing the figures as literal strings (matching the
# original inlineStr cells) instead of coercing them into numbers/percentages.
$updates = @{
    'D2' = '303.64'
    'E2' = '5.73%'
    'D3' = '32.00'
    'E3' = '8.80%'
    'D4' = '5.297'
    'E4' = '3.78%'
    'D5' = '0.07481'
    'E5' = '8.62%'
    'D6' = '7.793'
    'E6' = '5.47%'
    'D7' = '3.805'
    'E7' = '9.16%'
    'E8' = '5.52%'
    'D9' = '0.9198'
    'E9' = '1.84%'
    'D10' = '0.01766'
    'E10' = '2,614.82%'
    'D11' = '0.1701'
    'E11' = '6.85%'
    'D12' = '0.07666'
    'E12' = '10.59%'
    'D13' = '0.08075'
    'E13' = '6.05%'
    'D14' = '0.03053'
    'E14' = '4.58%'
    'D15' = '0.09885'
    'E15' = '10.01%'
    'D16' = '0.001511'
    'E16' = '-4.69%'
    'D17' = '0.04561'
    'E17' = '1.39%'
    'D18' = '0.006401'
    'E18' = '0.18%'
    'D19' = '3.484'
    'E19' = '0.43%'
    'D20' = '2.228'
    'E20' = '-0.12%'
    'D21' = '0.3333'
    'E21' = '3.54%'
    'D22' = '0.1343'
    'E22' = '1.46%'
    'D23' = '4.499'
    'E23' = '12.55%'
    'E24' = '4.38%'
    'D25' = '0.001218'
    'D26' = '0.004421'
    'E26' = '0.92%'
    'D27' = '0.0001402'
    'E27' = '20.43%'
    'D28' = '0.0001739'
    'E28' = '7.94%'
    'D40' = '0.04534'
    'E40' = '5.69%'
    'D41' = '0.007194'
    'E41' = '5.58%'
    'D42' = '0.1340'
    'E42' = '7.63%'
    'D43' = '0.002233'
    'E43' = '1.56%'
    'D44' = '0.01270'
    'E44' = '9.45%'
    'D45' = '0.00006133'
    'E45' = '7.10%'
    'D46' = '1.873'
    'E46' = '-2.95%'
    'D47' = '0.01298'
    'E47' = '-0.23%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
